# AQI_SCALE.xlsx edit script
# Adds RGB / Hex columns (J, K) to the "data" sheet, fixes a wrong RGB
# value (G3), and leaves everything else untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# ---- fix the "Moderate" row RGB value (G3: 0 -> 255) ----
$ws.Range("G3").Value = 255

# ---- header row ----
$ws.Range("J1").Value = "RGB"
$ws.Range("J1").WrapText = $true
$ws.Range("J1").VerticalAlignment = -4108   # xlVAlignCenter, matches header style used by A1:I1

$ws.Range("K1").Value = "Hex"
$ws.Range("K1").WrapText = $true
$ws.Range("K1").VerticalAlignment = -4108

# ---- RGB concat formulas ----
$ws.Range("J2").Formula = '=CONCAT(G2,",",H2,",",I2)'
$ws.Range("J3:J6").Formula = '=CONCAT(G3,",",H3,",",I3)'
$ws.Range("J7").Formula = '=CONCAT(G7,",",H7,",",I7)'

# ---- Hex color values ----
$ws.Range("K2").Value = "#009966"
$ws.Range("K3").Value = "#ffdd33"
$ws.Range("K4").Value = "#ff9933"
$ws.Range("K5").Value = "#cc0033"
$ws.Range("K6").Value = "#660099"
$ws.Range("K7").Value = "#7e0035"

# ---- view state: active cell K7, scrolled so column E is leftmost ----
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("K7").Select()
